# Update crypto price/volume data per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.696.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.55%  "

$ws.Range("D3").Value = "'1.868.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "'327.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("E7").Value = "  +0.91%  "

$ws.Range("E8").Value = "  +1.35%  "

$ws.Range("D9").Value = "'0.07903"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.58%  "

$ws.Range("D10").Value = "'0.9702"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("D11").Value = "'22.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("D12").Value = "'1.839.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.97%  "

$ws.Range("D13").Value = "'5.740"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").Value = "'6.935"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.23%  "

$ws.Range("D15").Value = "'0.06972"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("D16").Value = "'88.37"
$ws.Range("D16").Style = "Normal"

$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").Value = "'0.00001009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("D19").Value = "'16.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").Value = "'28.714.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.55%  "

$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").Value = "'11.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.33%  "

$ws.Range("D24").Value = "'2.113"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.00%  "

$ws.Range("D25").Value = "'2.086.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.89%  "

$ws.Range("D26").Value = "'153.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").Value = "'19.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("D28").Value = "'5.721"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.45%  "

$ws.Range("D29").Value = "'2.000"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.77%  "

$ws.Range("D30").Value = "'119.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.25%  "

$ws.Range("D31").Value = "'0.09369"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("D32").Value = "'0.9380"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("D33").Value = "'5.323"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").Value = "'1.347"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.59%  "

$ws.Range("D35").Value = "'3.361"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("E36").Value = "  -2.60%  "

$ws.Range("D37").Value = "'0.02128"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("D38").Value = "'1.146"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").Value = "'7.895"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.52%  "

$ws.Range("D40").Value = "'0.5661"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.70%  "

$ws.Range("D41").Value = "'9.965"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("D42").Value = "'0.1787"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("D43").Value = "'0.07237"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.94%  "

$ws.Range("D44").Value = "'11.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.83%  "

$ws.Range("D45").Value = "'0.5315"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.74%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'2.136"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.81%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.136"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.29%  "

$ws.Range("D48").Value = "'1.850"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("D49").Value = "'113.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.88%  "

$ws.Range("D50").Value = "'2.352"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.72%  "

$ws.Range("E51").Value = "  +0.36%  "
